{"js": "// The document is one date paragraph followed by a single 20x5 practice\n// table; `body.paragraphs` walks both in document order (date, then each\n// table row left-to-right), so we can replace every run's text positionally\n// with the corresponding answer from the regenerated worksheet.\nconst NEW_VALUES = [\n  \"2024-10-26 Saturday\", \"75-17=58\", \"19+25=44\", \"70-22=48\", \"8+78=86\", \"29+56=85\",\n  \"85-26=59\", \"43+49=92\", \"27+24=51\", \"16+56=72\", \"51-49=2\", \"71-33=38\",\n  \"19+62=81\", \"62-36=26\", \"16+68=84\", \"46+37=83\", \"45+16=61\", \"75-59=16\",\n  \"6+65=71\", \"87-79=8\", \"38+37=75\", \"39+18=57\", \"12-7=5\", \"64+7=71\",\n  \"52-33=19\", \"82-63=19\", \"2+79=81\", \"18+18=36\", \"37-8=29\", \"67+5=72\",\n  \"12+29=41\", \"9+12=21\", \"43+49=92\", \"88-59=29\", \"17+34=51\", \"21-8=13\",\n  \"84-27=57\", \"47-38=9\", \"91-13=78\", \"20-14=6\", \"27+67=94\", \"64-18=46\",\n  \"93-75=18\", \"43-36=7\", \"48-19=29\", \"23-16=7\", \"71-34=37\", \"39+47=86\",\n  \"10-9=1\", \"78-49=29\", \"95-6=89\", \"22+49=71\", \"88+9=97\", \"43+9=52\",\n  \"53+18=71\", \"12+19=31\", \"53-25=28\", \"62-7=55\", \"48+17=65\", \"38+49=87\",\n  \"96-27=69\", \"84-16=68\", \"86-49=37\", \"43-7=36\", \"83-25=58\", \"60-26=34\",\n  \"48+4=52\", \"63-18=45\", \"75+16=91\", \"46+48=94\", \"7+66=73\", \"12-4=8\",\n  \"16+79=95\", \"17+46=63\", \"58+9=67\", \"36+6=42\", \"46+19=65\", \"58-19=39\",\n  \"61-44=17\", \"28+37=65\", \"81-48=33\", \"71-68=3\", \"15-6=9\", \"19+67=86\",\n  \"27+45=72\", \"5+77=82\", \"12+59=71\", \"39+5=44\", \"29+15=44\", \"50-16=34\",\n  \"8+76=84\", \"49+2=51\", \"39+26=65\", \"55-9=46\", \"18+78=96\", \"23+69=92\",\n  \"60-51=9\", \"81-36=45\", \"87-29=58\", \"37+58=95\", \"25+36=61\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== NEW_VALUES.length) {\n  throw new Error(\n    \"Expected \" + NEW_VALUES.length + \" paragraphs but found \" +\n    paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(NEW_VALUES[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document is one date paragraph followed by a single 20x5 practice\n# table. Update the date line, then replace every cell's answer in place\n# (row-major order) using the regenerated worksheet's values.\n$d = $word.ActiveDocument\n\n# Update the date line (first paragraph, outside the table).\n$d.Paragraphs.Item(1).Range.Text = '2024-10-26 Saturday'\n\n# New answers for the 20x5 practice table, in row-major order.\n$newValues = @(\n    @('75-17=58', '19+25=44', '70-22=48', '8+78=86', '29+56=85'),\n    @('85-26=59', '43+49=92', '27+24=51', '16+56=72', '51-49=2'),\n    @('71-33=38', '19+62=81', '62-36=26', '16+68=84', '46+37=83'),\n    @('45+16=61', '75-59=16', '6+65=71', '87-79=8', '38+37=75'),\n    @('39+18=57', '12-7=5', '64+7=71', '52-33=19', '82-63=19'),\n    @('2+79=81', '18+18=36', '37-8=29', '67+5=72', '12+29=41'),\n    @('9+12=21', '43+49=92', '88-59=29', '17+34=51', '21-8=13'),\n    @('84-27=57', '47-38=9', '91-13=78', '20-14=6', '27+67=94'),\n    @('64-18=46', '93-75=18', '43-36=7', '48-19=29', '23-16=7'),\n    @('71-34=37', '39+47=86', '10-9=1', '78-49=29', '95-6=89'),\n    @('22+49=71', '88+9=97', '43+9=52', '53+18=71', '12+19=31'),\n    @('53-25=28', '62-7=55', '48+17=65', '38+49=87', '96-27=69'),\n    @('84-16=68', '86-49=37', '43-7=36', '83-25=58', '60-26=34'),\n    @('48+4=52', '63-18=45', '75+16=91', '46+48=94', '7+66=73'),\n    @('12-4=8', '16+79=95', '17+46=63', '58+9=67', '36+6=42'),\n    @('46+19=65', '58-19=39', '61-44=17', '28+37=65', '81-48=33'),\n    @('71-68=3', '15-6=9', '19+67=86', '27+45=72', '5+77=82'),\n    @('12+59=71', '39+5=44', '29+15=44', '50-16=34', '8+76=84'),\n    @('49+2=51', '39+26=65', '55-9=46', '18+78=96', '23+69=92'),\n    @('60-51=9', '81-36=45', '87-29=58', '37+58=95', '25+36=61')\n)\n\n$t = $d.Tables.Item(1)\nif ($t.Rows.Count -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) rows but found $($t.Rows.Count)\"\n}\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $row = $t.Rows.Item($r)\n    $rowValues = $newValues[$r - 1]\n    if ($row.Cells.Count -ne $rowValues.Count) {\n        throw \"Row $r`: expected $($rowValues.Count) cells but found $($row.Cells.Count)\"\n    }\n    for ($c = 1; $c -le $row.Cells.Count; $c++) {\n        $row.Cells.Item($c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
